# Weekly NYPD CompStat data refresh - "New crime data collected"
# Updates the report header (issue number + week-covering dates) and
# refreshes every Murder..Hate Crimes precinct-group figure for the
# current week, matching the latest weekly source pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text -------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Crime-complaint grid (rows 14-30, cols C-N) ------------------------
    $ws.Range("D14").Value = 2
    $ws.Range("G14").Value = 4
    $ws.Range("J14").Value = 12
    $ws.Range("K14").Value = -33.333333333333
    $ws.Range("L14").Value = -42.857142857142
    $ws.Range("N14").Value = -83.673469387755
    $ws.Range("C15").Value = 8
    $ws.Range("D15").Value = 1
    $ws.Range("E15").Value = 700
    $ws.Range("F15").Value = 15
    $ws.Range("G15").Value = 11
    $ws.Range("H15").Value = 36.363636363636
    $ws.Range("I15").Value = 76
    $ws.Range("J15").Value = 61
    $ws.Range("K15").Value = 24.590163934426
    $ws.Range("L15").Value = 22.580645161290
    $ws.Range("M15").Value = 80.952380952380
    $ws.Range("N15").Value = 18.75
    $ws.Range("C16").Value = 33
    $ws.Range("D16").Value = 36
    $ws.Range("E16").Value = -8.333333333333
    $ws.Range("F16").Value = 143
    $ws.Range("G16").Value = 107
    $ws.Range("H16").Value = 33.644859813084
    $ws.Range("I16").Value = 677
    $ws.Range("J16").Value = 562
    $ws.Range("K16").Value = 20.462633451957
    $ws.Range("L16").Value = 70.959595959596
    $ws.Range("M16").Value = 1.957831325301
    $ws.Range("N16").Value = -79.803102625298
    $ws.Range("C17").Value = 56
    $ws.Range("D17").Value = 55
    $ws.Range("E17").Value = 1.818181818181
    $ws.Range("F17").Value = 223
    $ws.Range("G17").Value = 197
    $ws.Range("H17").Value = 13.197969543147
    $ws.Range("I17").Value = 992
    $ws.Range("J17").Value = 835
    $ws.Range("K17").Value = 18.802395209580
    $ws.Range("L17").Value = 54.276827371695
    $ws.Range("M17").Value = 82.352941176470
    $ws.Range("N17").Value = 0.813008130081
    $ws.Range("C18").Value = 36
    $ws.Range("D18").Value = 34
    $ws.Range("E18").Value = 5.882352941176
    $ws.Range("F18").Value = 144
    $ws.Range("G18").Value = 143
    $ws.Range("H18").Value = 0.699300699300
    $ws.Range("I18").Value = 749
    $ws.Range("J18").Value = 702
    $ws.Range("K18").Value = 6.695156695156
    $ws.Range("L18").Value = 22.786885245901
    $ws.Range("M18").Value = -24.572004028197
    $ws.Range("N18").Value = -86.175710594315
    $ws.Range("C19").Value = 120
    $ws.Range("D19").Value = 120
    $ws.Range("E19").Value = 0
    $ws.Range("F19").Value = 497
    $ws.Range("G19").Value = 493
    $ws.Range("H19").Value = 0.811359026369
    $ws.Range("I19").Value = 2407
    $ws.Range("J19").Value = 2578
    $ws.Range("K19").Value = -6.633048875096
    $ws.Range("L19").Value = 73.790613718411
    $ws.Range("M19").Value = 70.467422096317
    $ws.Range("N19").Value = -13.912732474964
    $ws.Range("C20").Value = 54
    $ws.Range("E20").Value = 42.105263157894
    $ws.Range("F20").Value = 202
    $ws.Range("G20").Value = 138
    $ws.Range("H20").Value = 46.376811594202
    $ws.Range("I20").Value = 818
    $ws.Range("J20").Value = 553
    $ws.Range("K20").Value = 47.920433996383
    $ws.Range("L20").Value = 100.490196078431
    $ws.Range("M20").Value = 28.818897637795
    $ws.Range("N20").Value = -90.473972283684
    $ws.Range("C21").Value = 307
    $ws.Range("D21").Value = 286
    $ws.Range("E21").Value = 7.342657342657
    $ws.Range("F21").Value = 1224
    $ws.Range("G21").Value = 1093
    $ws.Range("H21").Value = 11.985361390667
    $ws.Range("I21").Value = 5727
    $ws.Range("J21").Value = 5303
    $ws.Range("K21").Value = 7.995474259852
    $ws.Range("L21").Value = 62.791358726549
    $ws.Range("M21").Value = 33.186046511627
    $ws.Range("N21").Value = -73.049411764705
    $ws.Range("C22").Value = 17
    $ws.Range("D22").Value = 2
    $ws.Range("E22").Value = 750
    $ws.Range("F22").Value = 32
    $ws.Range("G22").Value = 14
    $ws.Range("H22").Value = 128.571428571429
    $ws.Range("I22").Value = 134
    $ws.Range("J22").Value = 86
    $ws.Range("K22").Value = 55.813953488372
    $ws.Range("L22").Value = 211.627906976744
    $ws.Range("M22").Value = 86.111111111111
    $ws.Range("C23").Value = 7
    $ws.Range("D23").Value = 4
    $ws.Range("E23").Value = 75
    $ws.Range("F23").Value = 29
    $ws.Range("G23").Value = 17
    $ws.Range("H23").Value = 70.588235294117
    $ws.Range("I23").Value = 95
    $ws.Range("J23").Value = 90
    $ws.Range("K23").Value = 5.555555555555
    $ws.Range("L23").Value = 43.939393939393
    $ws.Range("M23").Value = 86.274509803921
    $ws.Range("C24").Value = 270
    $ws.Range("D24").Value = 320
    $ws.Range("E24").Value = -15.625
    $ws.Range("F24").Value = 1098
    $ws.Range("G24").Value = 1212
    $ws.Range("H24").Value = -9.405940594059
    $ws.Range("I24").Value = 5460
    $ws.Range("J24").Value = 5399
    $ws.Range("K24").Value = 1.129838859047
    $ws.Range("L24").Value = 34.981458590852
    $ws.Range("M24").Value = 66.921430755120
    $ws.Range("C25").Value = 108
    $ws.Range("D25").Value = 89
    $ws.Range("E25").Value = 21.348314606741
    $ws.Range("F25").Value = 447
    $ws.Range("G25").Value = 387
    $ws.Range("H25").Value = 15.503875968992
    $ws.Range("I25").Value = 1855
    $ws.Range("J25").Value = 1735
    $ws.Range("K25").Value = 6.916426512968
    $ws.Range("L25").Value = 39.159789947486
    $ws.Range("M25").Value = 8.226371061843
    $ws.Range("C26").Value = 10
    $ws.Range("D26").Value = 3
    $ws.Range("E26").Value = 233.333333333333
    $ws.Range("F26").Value = 26
    $ws.Range("G26").Value = 23
    $ws.Range("H26").Value = 13.043478260869
    $ws.Range("I26").Value = 116
    $ws.Range("J26").Value = 95
    $ws.Range("K26").Value = 22.105263157894
    $ws.Range("L26").Value = 20.833333333333
    $ws.Range("C27").Value = 27
    $ws.Range("D27").Value = 11
    $ws.Range("E27").Value = 145.454545454545
    $ws.Range("F27").Value = 73
    $ws.Range("G27").Value = 32
    $ws.Range("H27").Value = 128.125
    $ws.Range("I27").Value = 263
    $ws.Range("J27").Value = 203
    $ws.Range("K27").Value = 29.556650246305
    $ws.Range("L27").Value = 33.502538071066
    $ws.Range("C28").Value = 2
    $ws.Range("D28").Value = 3
    $ws.Range("E28").Value = -33.333333333333
    $ws.Range("I28").Value = 26
    $ws.Range("J28").Value = 25
    $ws.Range("K28").Value = 4
    $ws.Range("L28").Value = 36.842105263157
    $ws.Range("M28").Value = 62.5
    $ws.Range("N28").Value = -71.739130434782
    $ws.Range("C29").Value = 2
    $ws.Range("D29").Value = 2
    $ws.Range("G29").Value = 5
    $ws.Range("H29").Value = -20
    $ws.Range("I29").Value = 24
    $ws.Range("J29").Value = 23
    $ws.Range("K29").Value = 4.347826086956
    $ws.Range("L29").Value = 50
    $ws.Range("M29").Value = 100
    $ws.Range("N29").Value = -72.727272727272
    $ws.Range("C30").Value = "0"
    $ws.Range("D30").Value = "0"
    $ws.Range("E30").Value = "***.*"
    $ws.Range("F30").Value = 6
    $ws.Range("H30").Value = 200
    $ws.Range("I30").Value = 31
    $ws.Range("K30").Value = 29.166666666666
    $ws.Range("L30").Value = 10.714285714285

Write-Host "edit applied"
